$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 19.00506019592285
$ws.Range("C3").Value = 17.84777641296387
$ws.Range("C4").Value = 17.79007911682129
$ws.Range("C5").Value = 17.98081398010254
$ws.Range("C6").Value = 18.0962085723877
